$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add the two new data rows (Iphone 12 / Iphone 13), copying the existing
#    row format down first so the new rows pick up the same font/fill/border
#    as the rest of the table body.
# ---------------------------------------------------------------------------
$ws.Range("A16:D16").Copy()
$ws.Range("A18:D18").PasteSpecial(-4122)
$ws.Range("A16:D16").Copy()
$ws.Range("A19:D19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A18").Value = "Iphone 12"
$ws.Range("B18").Value = 2020
$ws.Range("C18").Value = 44857
$ws.Range("D18").Formula = "=DATE(B18, MONTH(C18), DAY(C18))"

$ws.Range("A19").Value = "Iphone 13"
$ws.Range("B19").Value = 2021
$ws.Range("C19").Value = 44828
$ws.Range("D19").Formula = "=DATE(B19, MONTH(C19), DAY(C19))"

# ---------------------------------------------------------------------------
# 2. Bump the font size for the whole table (header + body + date column) to
#    20pt, and make row 17 (previously styled with a darker bottom border)
#    consistent with the rest of the body rows.
# ---------------------------------------------------------------------------
$ws.Range("A1:D19").Font.Size = 20
$ws.Range("A17:C17").Borders.Item(9).Color = 15658734

# ---------------------------------------------------------------------------
# 3. Column widths / row heights to match the new, larger layout.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 19.9167
$ws.Columns.Item(2).ColumnWidth = 18.9167
$ws.Columns.Item(3).ColumnWidth = 23.0834
$ws.Columns.Item(4).ColumnWidth = 39.4167

$ws.Rows.Item(1).RowHeight = 27
$ws.Rows.Item(2).RowHeight = 27
$ws.Rows.Item(3).RowHeight = 27
$ws.Rows.Item(4).RowHeight = 27
$ws.Rows.Item(5).RowHeight = 27
$ws.Rows.Item(6).RowHeight = 27
$ws.Rows.Item(7).RowHeight = 27
$ws.Rows.Item(8).RowHeight = 27
$ws.Rows.Item(9).RowHeight = 27
$ws.Rows.Item(10).RowHeight = 27
$ws.Rows.Item(11).RowHeight = 27
$ws.Rows.Item(12).RowHeight = 27
$ws.Rows.Item(13).RowHeight = 27
$ws.Rows.Item(14).RowHeight = 27
$ws.Rows.Item(15).RowHeight = 27
$ws.Rows.Item(16).RowHeight = 27
$ws.Rows.Item(17).RowHeight = 27
$ws.Rows.Item(18).RowHeight = 27
$ws.Rows.Item(19).RowHeight = 27

# ---------------------------------------------------------------------------
# 4. Selection cosmetics to match the saved view.
# ---------------------------------------------------------------------------
$ws.Range("D18").Select()
